$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2..51: columns B (Coin), C (Link), D (Price), E (Volume(1h)), G (Hora)
# Values below are the post-edit (target) contents for each cell, taken from the
# GitHub-Actions "Updated symbol list" refresh. Every cell is written as literal
# text (leading apostrophe + Style reset keeps numeric-looking strings such as
# "246.49" or "17" stored as text instead of being auto-converted to numbers).

$rows = @(
    @{ Row = 2;  B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "246.49"; E = "1BNBBNB" },
    @{ Row = 3;  B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "24.21"; E = "2OKBOKB" },
    @{ Row = 4;  B = "HuobiToken"; C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D = "5.288"; E = "3HuobiTokenHT" },
    @{ Row = 5;  B = "Cronos"; C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D = "0.05816"; E = "4CronosCRO" },
    @{ Row = 6;  B = "KuCoinToken"; C = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"; D = "6.476"; E = "5KuCoinTokenKCS" },
    @{ Row = 7;  B = "GateToken"; C = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; D = "3.135"; E = "6GateTokenGT" },
    @{ Row = 8;  B = "MXToken"; C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D = "0.8182"; E = "7MXTokenMX" },
    @{ Row = 9;  B = "FTXToken"; C = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; D = "0.8776"; E = "8FTXTokenFTT" },
    @{ Row = 10; B = "One"; C = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; D = "0.01005"; E = "9OneONEBestin24h" },
    @{ Row = 11; B = "WazirX"; C = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; D = "0.1379"; E = "10WazirXWRX" },
    @{ Row = 12; B = "MandalaExchangeToken"; C = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D = "0.06996"; E = "11MandalaExchangeTokenMDX" },
    @{ Row = 13; B = "LiechtensteinCryptoassetsExchange"; C = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; D = "0.03122"; E = "12LiechtensteinCryptoassetsExchangeLCX" },
    @{ Row = 14; B = "BitrueCoin"; C = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D = "0.02944"; E = "13BitrueCoinBTR" },
    @{ Row = 15; B = "BitMartToken"; C = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D = "0.09409"; E = "14BitMartTokenBMX" },
    @{ Row = 16; B = "MCDex"; C = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"; D = "3.748"; E = "15MCDexMCB" },
    @{ Row = 17; B = "BitForexToken"; C = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D = "0.001526"; E = "16BitForexTokenBF" },
    @{ Row = 18; B = "CoinExToken"; C = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; D = "0.04664"; E = "17CoinExTokenCET" },
    @{ Row = 19; B = "TigerCash"; C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D = "0.006138"; E = "18TigerCashTCH" },
    @{ Row = 20; B = "BitKan"; C = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"; D = "0.001241"; E = "19BitKanKAN" },
    @{ Row = 21; B = "HotbitToken"; C = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"; D = "0.004669"; E = "20HotbitTokenHTB" },
    @{ Row = 22; B = "NitroEx"; C = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"; D = "0.00006104"; E = "21NitroExNTX" },
    @{ Row = 23; B = "LEO"; C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D = "3.531"; E = "22LEOLEO" },
    @{ Row = 24; B = "BTSEToken"; C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; D = "2.142"; E = "23BTSETokenBTSE" },
    @{ Row = 25; B = "BitpandaEcosystemToken"; C = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; D = "0.3181"; E = "24BitpandaEcosystemTokenBEST" },
    @{ Row = 26; B = "ProBitToken"; C = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"; D = "0.1306"; E = "25ProBitTokenPROB" },
    @{ Row = 27; B = "ZBToken"; C = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"; D = "0.1329"; E = "26ZBTokenZB" },
    @{ Row = 28; B = "UpBots"; C = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"; D = "0.0002333"; E = "27UpBotsUBXT" },
    @{ Row = 29; B = "Spectre.aiUtilityToken"; C = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"; D = "--"; E = "28Spectre.aiUtilityTokenSXUT" },
    @{ Row = 30; B = "LegolasExchange"; C = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"; D = "--"; E = "29LegolasExchangeLGO" },
    @{ Row = 31; B = "BitZToken"; C = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"; D = "--"; E = "30BitZTokenBZ" },
    @{ Row = 32; B = "Birake"; C = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"; D = "--"; E = "31BirakeBIR" },
    @{ Row = 33; B = "NashExchange"; C = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"; D = "--"; E = "32NashExchangeNEX" },
    @{ Row = 34; B = "AAXToken"; C = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"; D = "--"; E = "33AAXTokenAAB" },
    @{ Row = 35; B = "CenX"; C = "https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"; D = "--"; E = "34CenXCENX" },
    @{ Row = 36; B = "BNIXToken"; C = "https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"; D = "--"; E = "35BNIXTokenBNIX" },
    @{ Row = 37; B = "Polkally"; C = "https://coinranking.com/coin/NkDWUL8F-+polkally-kally"; D = "--"; E = "36PolkallyKALLY" },
    @{ Row = 38; B = "Charli3"; C = "https://coinranking.com/coin/8SgjMSqUk+charli3-c3"; D = "--"; E = "37Charli3C3" },
    @{ Row = 39; B = "BlubitexToken"; C = "https://coinranking.com/coin/Y9oImHIW5+blubitextoken-bbe"; D = "--"; E = "38BlubitexTokenBBE" },
    @{ Row = 40; B = "IDEX"; C = "https://coinranking.com/coin/ZiRElvGxqQaf+idex-idex"; D = "0.03727"; E = "39IDEXIDEX" },
    @{ Row = 41; B = "KickToken"; C = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"; D = "0.006364"; E = "40KickTokenKICK" },
    @{ Row = 42; B = "BKEXToken"; C = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"; D = "0.1057"; E = "41BKEXTokenBKK" },
    @{ Row = 43; B = "CEJI"; C = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"; D = "0.003001"; E = "42CEJICEJI" },
    @{ Row = 44; B = "LocalTraders"; C = "https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct"; D = "0.007869"; E = "43LocalTradersLCT" },
    @{ Row = 45; B = "CoinLion"; C = "https://coinranking.com/coin/sot4vgRyjNXek+coinlion-lion"; D = "0.00005259"; E = "44CoinLionLION" },
    @{ Row = 46; B = "Kangarootoken"; C = "https://coinranking.com/coin/zkVNkSGwZ3+kangarootoken-gar"; D = "0.00000000751"; E = "45KangarootokenGAR" },
    @{ Row = 47; B = "CoinbaseStockToken"; C = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"; D = "0.4053"; E = "46CoinbaseStockTokenCOIN" },
    @{ Row = 48; B = "BOLO"; C = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"; D = "0.002427"; E = "47BOLOBOLOWorstin24h" },
    @{ Row = 49; B = "CryptobidCoin"; C = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"; D = "0.00002101"; E = "48CryptobidCoinCBC" },
    @{ Row = 50; B = "SpecialPowerGold"; C = "https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg"; D = "0.0002001"; E = "49SpecialPowerGoldSPG" },
    @{ Row = 51; B = "DigiFinexToken"; C = "https://coinranking.com/coin/rY6dWXQL4+digifinextoken-dft"; D = "--"; E = "50DigiFinexTokenDFT" }
)

foreach ($item in $rows) {
    $r = $item.Row

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = "'" + $item.B
    $cellB.Style = "Normal"

    $cellC = $ws.Cells.Item($r, 3)
    $cellC.Value = "'" + $item.C
    $cellC.Style = "Normal"

    $cellD = $ws.Cells.Item($r, 4)
    $cellD.Value = "'" + $item.D
    $cellD.Style = "Normal"

    $cellE = $ws.Cells.Item($r, 5)
    $cellE.Value = "'" + $item.E
    $cellE.Style = "Normal"

    $cellG = $ws.Cells.Item($r, 7)
    $cellG.Value = "'17"
    $cellG.Style = "Normal"
}
